# Ngf-Ngfr NATMI LR-pairs sheet: refreshed with new TPM-based recalculation.
# The "ECs" target-cluster rows are dropped (Ligand-expressing / receptor
# stats were recomputed), shrinking the table from 9 data rows (3 senders x
# 3 targets) to 6 (3 senders x {FAPs, MuSCs} targets), and every remaining
# numeric column is refreshed with the new TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7 (columns A..T), taken from the recalculated NATMI output.
$data = @(
    @("ECs",   "Ngf", "Ngfr", "FAPs",  2, 0.6666666666666666, 0.1930666666666666, 0.5791999999999999, 0.01292026122037801, 0.01292026122037801, 3, 1, 1.311698333333333,  3.935095, 0.6472031793931536, 0.6472031793931535, 0.2532452248888888,  2.279207024,        0.008362034140418712, 0.008362034140418712),
    @("ECs",   "Ngf", "Ngfr", "MuSCs", 2, 0.6666666666666666, 0.1930666666666666, 0.5791999999999999, 0.01292026122037801, 0.01292026122037801, 3, 1, 0.7150196666666666, 2.145059, 0.3527968206068465, 0.3527968206068464, 0.1380464636444444,  1.2424181728,       0.004558227079959295, 0.004558227079959295),
    @("FAPs",  "Ngf", "Ngfr", "FAPs",  2, 0.6666666666666666, 0.666149,           1.998447,            0.04457951877603724, 0.04457951877603725, 3, 1, 1.311698333333333,  3.935095, 0.6472031793931536, 0.6472031793931535, 0.8737865330516668,  7.864078797465,     0.02885200628766809,  0.02885200628766809),
    @("FAPs",  "Ngf", "Ngfr", "MuSCs", 2, 0.6666666666666666, 0.666149,           1.998447,            0.04457951877603724, 0.04457951877603725, 3, 1, 0.7150196666666666, 2.145059, 0.3527968206068465, 0.3527968206068464, 0.4763096359303333,  4.286786723373,     0.01572751248836915,  0.01572751248836915),
    @("MuSCs", "Ngf", "Ngfr", "FAPs",  3, 1,                  14.08372266666667,  42.251168,           0.9425002200035847,  0.9425002200035848,  3, 1, 1.311698333333333,  3.935095, 0.6472031793931536, 0.6472031793931535, 18.47359554899556,   166.26235994096,    0.6099891389650668,   0.6099891389650668),
    @("MuSCs", "Ngf", "Ngfr", "MuSCs", 3, 1,                  14.08372266666667,  42.251168,           0.9425002200035847,  0.9425002200035848,  3, 1, 0.7150196666666666, 2.145059, 0.3527968206068465, 0.3527968206068464, 10.07013868654578,   90.63124817891199,  0.332511081038518,    0.332511081038518)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

# Drop the old rows 8-10 (previously the "MuSCs sending" x "ECs/FAPs/MuSCs
# target" block's tail); the table now ends at row 7.
$ws.Rows("8:10").Delete()
